$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.203.81'
$ws.Range("E2").Value = '  -5.34%  '
$ws.Range("D3").Value = '3.343.10'
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("D5").Value = '566.11'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = '131.86'
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.341.85'
$ws.Range("E8").Value = '  -2.52%  '
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").Value = '7.46'
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("D11").Value = '0.119'
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").Value = '3.908.95'
$ws.Range("E13").Value = '  -2.39%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '3.338.86'
$ws.Range("E15").Value = '  -2.54%  '
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("D17").Value = '24.87'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").Value = '60.265.66'
$ws.Range("E18").Value = '  -5.26%  '
$ws.Range("D19").Value = '13.59'
$ws.Range("E19").Value = '  +2.52%  '
$ws.Range("D20").Value = '5.71'
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("D21").Value = '9.26'
$ws.Range("E21").Value = '  -5.73%  '
$ws.Range("D22").Value = '355.63'
$ws.Range("E22").Value = '  -7.36%  '
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D25").Value = '3.469.10'
$ws.Range("E25").Value = '  -2.67%  '
$ws.Range("D26").Value = '69.43'
$ws.Range("E26").Value = '  -5.94%  '
$ws.Range("D27").Value = '0.0000112'
$ws.Range("E27").Value = '  +2.61%  '
$ws.Range("E28").Value = '  +17.05%  '
$ws.Range("E29").Value = '  +6.68%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").Value = '7.99'
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").Value = '3.369.81'
$ws.Range("E35").Value = '  -2.53%  '
$ws.Range("D36").Value = '22.97'
$ws.Range("E36").Value = '  +1.37%  '
$ws.Range("E37").Value = '  +3.61%  '
$ws.Range("E38").Value = '  +2.38%  '
$ws.Range("E39").Value = '  +0.82%  '
$ws.Range("D40").Value = '159.04'
$ws.Range("D41").Value = '0.0774'
$ws.Range("E41").Value = '  +1.48%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '4.39'
$ws.Range("E43").Value = '  +2.04%  '
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Value = '1.19'
$ws.Range("E44").Value = '  +8.52%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '40.94'
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("D46").Value = '0.749'
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("D47").Value = '23.61'
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("D48").Value = '1.59'
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("E49").Value = '  +1.76%  '
$ws.Range("D50").Value = '22.57'
$ws.Range("E50").Value = '  +11.60%  '
$ws.Range("D51").Value = '0.895'
$ws.Range("E51").Value = '  +1.47%  '
